{"js": "const pairs = [\n  [\"2024-04-12 Friday\", \"2024-04-13 Saturday\"],\n  [\"29+30=\", \"61-11=\"],\n  [\"29+67=\", \"24+9=\"],\n  [\"17-4=\", \"5+86=\"],\n  [\"5+32=\", \"86-2=\"],\n  [\"46+14=\", \"72-25=\"],\n  [\"71+2=\", \"31+1=\"],\n  [\"44-31=\", \"64+4=\"],\n  [\"22+2=\", \"64+28=\"],\n  [\"49+1=\", \"35+45=\"],\n  [\"18+24=\", \"13-6=\"],\n  [\"61+17=\", \"70+21=\"],\n  [\"98-0=\", \"17-2=\"],\n  [\"17+3=\", \"29-17=\"],\n  [\"25-2=\", \"20+5=\"],\n  [\"24-20=\", \"37+35=\"],\n  [\"63-57=\", \"95-22=\"],\n  [\"78+14=\", \"94-62=\"],\n  [\"98-59=\", \"92-47=\"],\n  [\"41+28=\", \"50+44=\"],\n  [\"70-25=\", \"35-20=\"],\n  [\"33+11=\", \"34+40=\"],\n  [\"31+20=\", \"3+73=\"],\n  [\"67+23=\", \"31+24=\"],\n  [\"11+45=\", \"61+29=\"],\n  [\"78-69=\", \"67-58=\"],\n  [\"98-12=\", \"99-16=\"],\n  [\"72-11=\", \"50-31=\"],\n  [\"39-15=\", \"25+71=\"],\n  [\"5+8=\", \"42-17=\"],\n  [\"89-54=\", \"82-26=\"],\n  [\"83-62=\", \"42+5=\"],\n  [\"19+7=\", \"83-24=\"],\n  [\"68-37=\", \"2+61=\"],\n  [\"5+92=\", \"16+37=\"],\n  [\"62-19=\", \"79-74=\"],\n  [\"78-16=\", \"73-47=\"],\n  [\"65-19=\", \"82-39=\"],\n  [\"41+30=\", \"93-76=\"],\n  [\"65-24=\", \"89-61=\"],\n  [\"22+75=\", \"37+62=\"],\n  [\"26+50=\", \"99-26=\"],\n  [\"12-3=\", \"93+4=\"],\n  [\"81+5=\", \"65-52=\"],\n  [\"43-14=\", \"10+19=\"],\n  [\"8+17=\", \"59+0=\"],\n  [\"30+40=\", \"26+68=\"],\n  [\"98-54=\", \"24+70=\"],\n  [\"80-18=\", \"59+34=\"],\n  [\"62+28=\", \"56-52=\"],\n  [\"85-27=\", \"28+1=\"],\n  [\"98-20=\", \"22-22=\"],\n  [\"2+90=\", \"60+10=\"],\n  [\"68-42=\", \"90-5=\"],\n  [\"96-41=\", \"60+26=\"],\n  [\"42-12=\", \"87-54=\"],\n  [\"11+18=\", \"42+37=\"],\n  [\"77-56=\", \"6+57=\"],\n  [\"41-35=\", \"63-23=\"],\n  [\"19-4=\", \"41-16=\"],\n  [\"42+53=\", \"69-11=\"],\n  [\"34-28=\", \"48+38=\"],\n  [\"67+31=\", \"67+22=\"],\n  [\"45-15=\", \"80-71=\"],\n  [\"37+55=\", \"55+29=\"],\n  [\"7+62=\", \"77+9=\"],\n  [\"21-20=\", \"23-8=\"],\n  [\"41-40=\", \"89+2=\"],\n  [\"88-48=\", \"26+21=\"],\n  [\"63-36=\", \"96-8=\"],\n  [\"0+29=\", \"20+68=\"],\n  [\"84-24=\", \"64+20=\"],\n  [\"77-22=\", \"15+19=\"],\n  [\"52+14=\", \"43+55=\"],\n  [\"34-25=\", \"93-34=\"],\n  [\"81+17=\", \"80-38=\"],\n  [\"39-33=\", \"86-82=\"],\n  [\"51+36=\", \"78-13=\"],\n  [\"0+22=\", \"83-24=\"],\n  [\"84+0=\", \"63+34=\"],\n  [\"14+11=\", \"49+45=\"],\n  [\"58+7=\", \"45-33=\"],\n  [\"59-39=\", \"79+0=\"],\n  [\"12+59=\", \"36+59=\"],\n  [\"78+11=\", \"10+14=\"],\n  [\"4+21=\", \"50-21=\"],\n  [\"46+33=\", \"60+28=\"],\n  [\"5+5=\", \"42+37=\"],\n  [\"15-8=\", \"16-15=\"],\n  [\"78-45=\", \"13+77=\"],\n  [\"52-21=\", \"55-40=\"],\n  [\"8+3=\", \"1+10=\"],\n  [\"66-42=\", \"18+45=\"],\n  [\"19+8=\", \"94-56=\"],\n  [\"51-10=\", \"62+22=\"],\n  [\"25+61=\", \"57+31=\"],\n  [\"34+43=\", \"68-43=\"],\n  [\"22+47=\", \"66+1=\"],\n  [\"90-47=\", \"69-37=\"],\n  [\"48+47=\", \"90-55=\"],\n  [\"86+0=\", \"89-72=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: locate every target range using the ORIGINAL (pre-edit) text.\n// Doing all the searches before any text is changed avoids a later\n// replacement's new text accidentally containing an earlier search target\n// as a substring (e.g. '7+62=' is a substring of '37+62=').\nconst searchResults = [];\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  searchResults.push(results);\n}\nawait context.sync();\n\n// Phase 2: replace each located range with its new text.\nfor (let i = 0; i < pairs.length; i++) {\n  const [before, after] = pairs[i];\n  const results = searchResults[i];\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${before}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# --- Update the date line (first paragraph) ---\n$dateExpected = \"2024-04-12 Friday\"\n$dateNew = \"2024-04-13 Saturday\"\n$p1 = $d.Paragraphs.Item(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13, [char]7)\nif ($p1Text -ne $dateExpected) {\n    throw \"Unexpected first paragraph text: $p1Text\"\n}\n$p1.Range.Text = $dateNew\n\n# --- Update the 100 arithmetic-problem table cells (20 rows x 5 columns, row-major) ---\n$values = @(\n    @(\"29+30=\", \"61-11=\"),\n    @(\"29+67=\", \"24+9=\"),\n    @(\"17-4=\", \"5+86=\"),\n    @(\"5+32=\", \"86-2=\"),\n    @(\"46+14=\", \"72-25=\"),\n    @(\"71+2=\", \"31+1=\"),\n    @(\"44-31=\", \"64+4=\"),\n    @(\"22+2=\", \"64+28=\"),\n    @(\"49+1=\", \"35+45=\"),\n    @(\"18+24=\", \"13-6=\"),\n    @(\"61+17=\", \"70+21=\"),\n    @(\"98-0=\", \"17-2=\"),\n    @(\"17+3=\", \"29-17=\"),\n    @(\"25-2=\", \"20+5=\"),\n    @(\"24-20=\", \"37+35=\"),\n    @(\"63-57=\", \"95-22=\"),\n    @(\"78+14=\", \"94-62=\"),\n    @(\"98-59=\", \"92-47=\"),\n    @(\"41+28=\", \"50+44=\"),\n    @(\"70-25=\", \"35-20=\"),\n    @(\"33+11=\", \"34+40=\"),\n    @(\"31+20=\", \"3+73=\"),\n    @(\"67+23=\", \"31+24=\"),\n    @(\"11+45=\", \"61+29=\"),\n    @(\"78-69=\", \"67-58=\"),\n    @(\"98-12=\", \"99-16=\"),\n    @(\"72-11=\", \"50-31=\"),\n    @(\"39-15=\", \"25+71=\"),\n    @(\"5+8=\", \"42-17=\"),\n    @(\"89-54=\", \"82-26=\"),\n    @(\"83-62=\", \"42+5=\"),\n    @(\"19+7=\", \"83-24=\"),\n    @(\"68-37=\", \"2+61=\"),\n    @(\"5+92=\", \"16+37=\"),\n    @(\"62-19=\", \"79-74=\"),\n    @(\"78-16=\", \"73-47=\"),\n    @(\"65-19=\", \"82-39=\"),\n    @(\"41+30=\", \"93-76=\"),\n    @(\"65-24=\", \"89-61=\"),\n    @(\"22+75=\", \"37+62=\"),\n    @(\"26+50=\", \"99-26=\"),\n    @(\"12-3=\", \"93+4=\"),\n    @(\"81+5=\", \"65-52=\"),\n    @(\"43-14=\", \"10+19=\"),\n    @(\"8+17=\", \"59+0=\"),\n    @(\"30+40=\", \"26+68=\"),\n    @(\"98-54=\", \"24+70=\"),\n    @(\"80-18=\", \"59+34=\"),\n    @(\"62+28=\", \"56-52=\"),\n    @(\"85-27=\", \"28+1=\"),\n    @(\"98-20=\", \"22-22=\"),\n    @(\"2+90=\", \"60+10=\"),\n    @(\"68-42=\", \"90-5=\"),\n    @(\"96-41=\", \"60+26=\"),\n    @(\"42-12=\", \"87-54=\"),\n    @(\"11+18=\", \"42+37=\"),\n    @(\"77-56=\", \"6+57=\"),\n    @(\"41-35=\", \"63-23=\"),\n    @(\"19-4=\", \"41-16=\"),\n    @(\"42+53=\", \"69-11=\"),\n    @(\"34-28=\", \"48+38=\"),\n    @(\"67+31=\", \"67+22=\"),\n    @(\"45-15=\", \"80-71=\"),\n    @(\"37+55=\", \"55+29=\"),\n    @(\"7+62=\", \"77+9=\"),\n    @(\"21-20=\", \"23-8=\"),\n    @(\"41-40=\", \"89+2=\"),\n    @(\"88-48=\", \"26+21=\"),\n    @(\"63-36=\", \"96-8=\"),\n    @(\"0+29=\", \"20+68=\"),\n    @(\"84-24=\", \"64+20=\"),\n    @(\"77-22=\", \"15+19=\"),\n    @(\"52+14=\", \"43+55=\"),\n    @(\"34-25=\", \"93-34=\"),\n    @(\"81+17=\", \"80-38=\"),\n    @(\"39-33=\", \"86-82=\"),\n    @(\"51+36=\", \"78-13=\"),\n    @(\"0+22=\", \"83-24=\"),\n    @(\"84+0=\", \"63+34=\"),\n    @(\"14+11=\", \"49+45=\"),\n    @(\"58+7=\", \"45-33=\"),\n    @(\"59-39=\", \"79+0=\"),\n    @(\"12+59=\", \"36+59=\"),\n    @(\"78+11=\", \"10+14=\"),\n    @(\"4+21=\", \"50-21=\"),\n    @(\"46+33=\", \"60+28=\"),\n    @(\"5+5=\", \"42+37=\"),\n    @(\"15-8=\", \"16-15=\"),\n    @(\"78-45=\", \"13+77=\"),\n    @(\"52-21=\", \"55-40=\"),\n    @(\"8+3=\", \"1+10=\"),\n    @(\"66-42=\", \"18+45=\"),\n    @(\"19+8=\", \"94-56=\"),\n    @(\"51-10=\", \"62+22=\"),\n    @(\"25+61=\", \"57+31=\"),\n    @(\"34+43=\", \"68-43=\"),\n    @(\"22+47=\", \"66+1=\"),\n    @(\"90-47=\", \"69-37=\"),\n    @(\"48+47=\", \"90-55=\"),\n    @(\"86+0=\", \"89-72=\")\n)\n\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $expected = $values[$idx][0]\n        $new = $values[$idx][1]\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -ne $expected) {\n            throw \"Cell R${r}C${c}: expected [$expected] but found [$cellText]\"\n        }\n        $cell.Range.Text = $new\n        $idx++\n    }\n}\n\nWrite-Output \"Replaced $idx table cells plus the date line.\"\n"}
